# Reorder the "Recorded By" (column G) author lists so that whenever the
# list currently starts with "System"/"system", the last author in the
# list is rotated to the front (e.g. "System, user@x.com" becomes
# "user@x.com, System"; "system, System, user@x.com" becomes
# "user@x.com, system, System"). Rows whose list does not start with
# "System" (or that only contain a single entry) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $rawParts = $text.Split(",")
    $parts = @()
    for ($i = 0; $i -lt $rawParts.Length; $i++) {
        $parts += $rawParts[$i].Trim()
    }

    if ($parts.Length -gt 1 -and $parts[0].ToLower() -eq "system") {
        $lastItem = $parts[$parts.Length - 1]
        $remaining = @()
        for ($i = 0; $i -lt $parts.Length - 1; $i++) {
            $remaining += $parts[$i]
        }
        $newParts = @($lastItem) + $remaining
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
